$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1042.4546
$ws.Range("I33").Value = 1042.4546
$ws.Range("K33").Value = 1042.4546
$ws.Range("M33").Value = -813.4546
$ws.Range("H97").Value = 1939.6
$ws.Range("I97").Value = 1500
$ws.Range("J97").Value = 2049.5
$ws.Range("K97").Value = 4500
$ws.Range("L97").Value = 6148.5
$ws.Range("M97").Value = -4004
$ws.Range("N97").Value = -7140.5
$ws.Range("H116").Value = 3833.3333
$ws.Range("J116").Value = 4250
$ws.Range("L116").Value = 4250
$ws.Range("N116").Value = -11134
$ws.Range("H135").Value = 3982
$ws.Range("I135").Value = 3982
$ws.Range("K135").Value = 35838
$ws.Range("M135").Value = -33303

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 149
$ws.Range("I4").Value = 190.33333
$ws.Range("J4").Value = 25
$ws.Range("K4").Value = 190.33333
$ws.Range("L4").Value = 25
$ws.Range("M4").Value = -74.33332999999999
$ws.Range("N4").Value = -257
$ws.Range("H45").Value = 13000
$ws.Range("J45").Value = 13000
$ws.Range("L45").Value = 13000
$ws.Range("N45").Value = -13754
$ws.Range("H112").Value = 31193
$ws.Range("J112").Value = 31193
$ws.Range("L112").Value = 31193
$ws.Range("N112").Value = -34147
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 500
$ws.Range("I4").Value = 500
$ws.Range("K4").Value = 500
$ws.Range("M4").Value = -388
$ws.Range("H5").Value = 199
$ws.Range("I5").Value = 165.33333
$ws.Range("J5").Value = 300
$ws.Range("K5").Value = 165.33333
$ws.Range("L5").Value = 300
$ws.Range("M5").Value = -53.33332999999999
$ws.Range("N5").Value = -524
$ws.Range("H7").Value = 128.42857
$ws.Range("I7").Value = 38.88889
$ws.Range("J7").Value = 289.6
$ws.Range("K7").Value = 38.88889
$ws.Range("L7").Value = 289.6
$ws.Range("M7").Value = 74.11111
$ws.Range("N7").Value = -515.6
$ws.Range("H31").Value = 8175.647
$ws.Range("I31").Value = 5333.3335
$ws.Range("J31").Value = 11373.25
$ws.Range("K31").Value = 5333.3335
$ws.Range("L31").Value = 11373.25
$ws.Range("M31").Value = -5038.3335
$ws.Range("N31").Value = -11963.25
$ws.Range("H34").Value = 8175.647
$ws.Range("I34").Value = 5333.3335
$ws.Range("J34").Value = 11373.25
$ws.Range("K34").Value = 5333.3335
$ws.Range("L34").Value = 11373.25
$ws.Range("M34").Value = -5131.3335
$ws.Range("N34").Value = -11777.25
$ws.Range("H53").Value = 26684
$ws.Range("J53").Value = 26684
$ws.Range("L53").Value = 26684
$ws.Range("N53").Value = -27898
$ws.Range("H58").Value = 5501.5386
$ws.Range("I58").Value = 1315.125
$ws.Range("K58").Value = 1315.125
$ws.Range("M58").Value = -1112.125
$ws.Range("H134").Value = 12110.667
$ws.Range("J134").Value = 10166.5
$ws.Range("L134").Value = 30499.5
$ws.Range("N134").Value = -35569.5
$ws.Range("H136").Value = 5501.5386
$ws.Range("I136").Value = 1315.125
$ws.Range("K136").Value = 3945.375
$ws.Range("M136").Value = -1395.375
$ws.Range("H141").Value = 279999
$ws.Range("J141").Value = 799999
$ws.Range("L141").Value = 799999
$ws.Range("N141").Value = -810359

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 200000430
$ws.Range("I4").Value = 695
$ws.Range("K4").Value = 2085
$ws.Range("M4").Value = -1973
$ws.Range("H7").Value = 22.75
$ws.Range("I7").Value = 19
$ws.Range("K7").Value = 57
$ws.Range("M7").Value = 55

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 36.066666
$ws.Range("I2").Value = 26.272728
$ws.Range("J2").Value = 63
$ws.Range("K2").Value = 26.272728
$ws.Range("L2").Value = 63
$ws.Range("M2").Value = 86.727272
$ws.Range("N2").Value = -289
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("N66").ClearContents()
$ws.Range("H131").Value = 70000
$ws.Range("J131").Value = 70000
$ws.Range("L131").Value = 70000
$ws.Range("N131").Value = -80080

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2676.111
$ws.Range("I22").Value = 2372.8572
$ws.Range("K22").Value = 2372.8572
$ws.Range("M22").Value = -2077.8572
$ws.Range("H27").Value = 2676.111
$ws.Range("I27").Value = 2372.8572
$ws.Range("K27").Value = 2372.8572
$ws.Range("M27").Value = -2265.8572
$ws.Range("H55").Value = 1257.5625
$ws.Range("I55").Value = 1141.4
$ws.Range("K55").Value = 1141.4
$ws.Range("M55").Value = -968.4000000000001
$ws.Range("H61").Value = 6835.3335
$ws.Range("I61").Value = 6835.3335
$ws.Range("K61").Value = 6835.3335
$ws.Range("M61").Value = -6633.3335
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H68").Value = 2000
$ws.Range("J68").Value = 2000
$ws.Range("L68").Value = 2000
$ws.Range("N68").Value = -3498
$ws.Range("H71").Value = 2000
$ws.Range("J71").Value = 2000
$ws.Range("L71").Value = 10000
$ws.Range("N71").Value = -17488
$ws.Range("H93").Value = 5331.6665
$ws.Range("I93").Value = 5331.6665
$ws.Range("K93").Value = 5331.6665
$ws.Range("M93").Value = -4083.6665
$ws.Range("H113").Value = 6835.3335
$ws.Range("I113").Value = 6835.3335
$ws.Range("K113").Value = 6835.3335
$ws.Range("M113").Value = -4665.3335
$ws.Range("H132").Value = 6982.591
$ws.Range("I132").Value = 6416.684
$ws.Range("K132").Value = 19250.052
$ws.Range("M132").Value = -16720.052
